$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 9292.134
$ws.Range("J112").Value = 12039.272
$ws.Range("L112").Value = 36117.81600000001
$ws.Range("N112").Value = -38333.81600000001
$ws.Range("H132").Value = 54753.18
$ws.Range("I132").Value = 67193.086
$ws.Range("K132").Value = 201579.258
$ws.Range("M132").Value = -199049.258
$ws.Range("H137").Value = 2005.0408
$ws.Range("I137").Value = 1690.9166
$ws.Range("K137").Value = 5072.7498
$ws.Range("M137").Value = -2522.7498
$ws.Range("H138").Value = 1663.8478
$ws.Range("I138").Value = 1085
$ws.Range("J138").Value = 3303.9167
$ws.Range("K138").Value = 3255
$ws.Range("L138").Value = 9911.750100000001
$ws.Range("M138").Value = 1885
$ws.Range("N138").Value = -20191.7501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9617167
$ws.Range("I32").Value = 10205538
$ws.Range("K32").Value = 10205538
$ws.Range("M32").Value = -10205251
$ws.Range("H61").Value = 10496.3125
$ws.Range("I61").Value = 22030
$ws.Range("K61").Value = 22030
$ws.Range("M61").Value = -21818
$ws.Range("H74").Value = 6801.0386
$ws.Range("I74").Value = 7030.4736
$ws.Range("J74").Value = 6178.2856
$ws.Range("K74").Value = 7030.4736
$ws.Range("L74").Value = 6178.2856
$ws.Range("M74").Value = -6156.4736
$ws.Range("N74").Value = -7926.2856
$ws.Range("H77").Value = 6801.0386
$ws.Range("I77").Value = 7030.4736
$ws.Range("J77").Value = 6178.2856
$ws.Range("K77").Value = 35152.368
$ws.Range("L77").Value = 30891.428
$ws.Range("M77").Value = -30784.368
$ws.Range("N77").Value = -39627.428
$ws.Range("H102").Value = 13756.03
$ws.Range("I102").Value = 14437.097
$ws.Range("J102").Value = 3199.5
$ws.Range("K102").Value = 14437.097
$ws.Range("L102").Value = 3199.5
$ws.Range("M102").Value = -12815.097
$ws.Range("N102").Value = -6443.5
$ws.Range("H132").Value = 2194.9883
$ws.Range("I132").Value = 2232.3247
$ws.Range("J132").Value = 1835.625
$ws.Range("K132").Value = 6696.9741
$ws.Range("L132").Value = 5506.875
$ws.Range("M132").Value = -4166.9741
$ws.Range("N132").Value = -10566.875
$ws.Range("H136").Value = 10496.3125
$ws.Range("I136").Value = 22030
$ws.Range("K136").Value = 66090
$ws.Range("M136").Value = -63540

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3020
$ws.Range("I86").Value = 2750
$ws.Range("J86").Value = 3200
$ws.Range("K86").Value = 2750
$ws.Range("L86").Value = 3200
$ws.Range("M86").Value = -1627
$ws.Range("N86").Value = -5446
$ws.Range("H89").Value = 3020
$ws.Range("I89").Value = 2750
$ws.Range("J89").Value = 3200
$ws.Range("K89").Value = 13750
$ws.Range("L89").Value = 16000
$ws.Range("M89").Value = -8134
$ws.Range("N89").Value = -27232
$ws.Range("H94").Value = 1560.1111
$ws.Range("I94").Value = 1518.2222
$ws.Range("J94").Value = 1643.8889
$ws.Range("K94").Value = 1518.2222
$ws.Range("L94").Value = 1643.8889
$ws.Range("M94").Value = -1067.2222
$ws.Range("N94").Value = -2545.8889
$ws.Range("H99").Value = 56002.273
$ws.Range("I99").Value = 85002.5
$ws.Range("J99").Value = 21202
$ws.Range("K99").Value = 85002.5
$ws.Range("L99").Value = 21202
$ws.Range("M99").Value = -83504.5
$ws.Range("N99").Value = -24198
$ws.Range("H107").Value = 22225482
$ws.Range("I107").Value = 3624.625
$ws.Range("K107").Value = 3624.625
$ws.Range("M107").Value = -1704.625
$ws.Range("H134").Value = 3073.7742
$ws.Range("I134").Value = 1390.8182
$ws.Range("J134").Value = 7187.6665
$ws.Range("K134").Value = 4172.4546
$ws.Range("L134").Value = 21562.9995
$ws.Range("M134").Value = -1637.4546
$ws.Range("N134").Value = -26632.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 70894.87
$ws.Range("I31").Value = 105439.35
$ws.Range("K31").Value = 105439.35
$ws.Range("M31").Value = -105144.35
$ws.Range("H34").Value = 70894.87
$ws.Range("I34").Value = 105439.35
$ws.Range("K34").Value = 105439.35
$ws.Range("M34").Value = -105237.35
$ws.Range("H36").Value = 12998.5
$ws.Range("I36").Value = 998
$ws.Range("K36").Value = 998
$ws.Range("M36").Value = -610
$ws.Range("H40").Value = 12998.5
$ws.Range("I40").Value = 998
$ws.Range("K40").Value = 998
$ws.Range("M40").Value = -838
$ws.Range("H134").Value = 13630.464
$ws.Range("I134").Value = 20373.234
$ws.Range("K134").Value = 61119.702
$ws.Range("M134").Value = -58584.702

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 188.8125
$ws.Range("J17").Value = 586.75
$ws.Range("L17").Value = 1760.25
$ws.Range("N17").Value = -2098.25
$ws.Range("H109").Value = 3067.3635
$ws.Range("I109").Value = 1975.1666
$ws.Range("K109").Value = 5925.4998
$ws.Range("M109").Value = -4885.4998
$ws.Range("H114").Value = 2786.1428
$ws.Range("I114").Value = 562.5
$ws.Range("J114").Value = 3675.6
$ws.Range("K114").Value = 1687.5
$ws.Range("L114").Value = 11026.8
$ws.Range("M114").Value = 1566.5
$ws.Range("N114").Value = -17534.8
$ws.Range("H117").Value = 4057.6667
$ws.Range("J117").Value = 5514.231
$ws.Range("L117").Value = 16542.693
$ws.Range("N117").Value = -23426.693
$ws.Range("H119").Value = 3337.25
$ws.Range("I119").Value = 3099.7144
$ws.Range("K119").Value = 9299.143199999999
$ws.Range("M119").Value = -4461.143199999999
$ws.Range("H132").Value = 3668.077
$ws.Range("I132").Value = 2882.8333
$ws.Range("J132").Value = 4341.143
$ws.Range("K132").Value = 25945.4997
$ws.Range("L132").Value = 39070.287
$ws.Range("M132").Value = -23415.4997
$ws.Range("N132").Value = -44130.287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2960.8667
$ws.Range("I80").Value = 3261.8572
$ws.Range("J80").Value = 2697.5
$ws.Range("K80").Value = 3261.8572
$ws.Range("L80").Value = 2697.5
$ws.Range("M80").Value = -2263.8572
$ws.Range("N80").Value = -4693.5
$ws.Range("H83").Value = 2960.8667
$ws.Range("I83").Value = 3261.8572
$ws.Range("J83").Value = 2697.5
$ws.Range("K83").Value = 16309.286
$ws.Range("L83").Value = 13487.5
$ws.Range("M83").Value = -11317.286
$ws.Range("N83").Value = -23471.5
$ws.Range("H97").Value = 2105.6667
$ws.Range("I97").Value = 1018.8947
$ws.Range("J97").Value = 6235.4
$ws.Range("K97").Value = 1018.8947
$ws.Range("L97").Value = 6235.4
$ws.Range("M97").Value = -522.8946999999999
$ws.Range("N97").Value = -7227.4
$ws.Range("H113").Value = 2762.6667
$ws.Range("I113").Value = 2401.1538
$ws.Range("J113").Value = 3702.6
$ws.Range("K113").Value = 2401.1538
$ws.Range("L113").Value = 3702.6
$ws.Range("M113").Value = -231.1538
$ws.Range("N113").Value = -8042.6
$ws.Range("H132").Value = 223714.16
$ws.Range("I132").Value = 305778.97
$ws.Range("J132").Value = 18552.125
$ws.Range("K132").Value = 917336.9099999999
$ws.Range("L132").Value = 55656.375
$ws.Range("M132").Value = -914806.9099999999
$ws.Range("N132").Value = -60716.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1000.5
$ws.Range("I46").Value = 1000.5
$ws.Range("K46").Value = 1000.5
$ws.Range("M46").Value = -812.5
$ws.Range("H93").Value = 2834.4443
$ws.Range("I93").Value = 2784.25
$ws.Range("J93").Value = 2934.8333
$ws.Range("K93").Value = 2784.25
$ws.Range("L93").Value = 2934.8333
$ws.Range("M93").Value = -1536.25
$ws.Range("N93").Value = -5430.8333
$ws.Range("H110").Value = 81596
$ws.Range("J110").Value = 81596
$ws.Range("L110").Value = 81596
$ws.Range("N110").Value = -89776
$ws.Range("H132").Value = 788163.8
$ws.Range("I132").Value = 1082758.8
$ws.Range("J132").Value = 2577.3333
$ws.Range("K132").Value = 3248276.4
$ws.Range("L132").Value = 7731.999899999999
$ws.Range("M132").Value = -3245746.4
$ws.Range("N132").Value = -12791.9999
$ws.Range("H136").Value = 32083.785
$ws.Range("I136").Value = 2273.3684
$ws.Range("K136").Value = 6820.1052
$ws.Range("M136").Value = -4270.1052

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3231.5
$ws.Range("I107").Value = 2026.7222
$ws.Range("K107").Value = 6080.1666
$ws.Range("M107").Value = -4160.1666
$ws.Range("H136").Value = 38067468
$ws.Range("I136").Value = 47497996
$ws.Range("J136").Value = 345344.5
$ws.Range("K136").Value = 142493988
$ws.Range("L136").Value = 1036033.5
$ws.Range("M136").Value = -142491438
$ws.Range("N136").Value = -1041133.5

$wb.Save()